$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (new quarters), shifting old D:K data to F:M
$ws.Columns("D:E").Insert()

# Copy number formats from column F (old column D) down into new D:E columns,
# restricted to the three data blocks (Income Statement / Balance Sheet / Cash Flow)
$ws.Range("F7:F35").Copy()
$ws.Range("D7:E35").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("D38:E77").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Give the two new columns a sensible width matching their neighbours
$ws.Range("D1:E1").ColumnWidth = $ws.Range("F1").ColumnWidth

# Populate the two new quarter columns (D = 2018-11-30, E = 2018-08-31) with data
$newData = @{
    7 = @(43434, 43343)
    8 = @(91500, 85300)
    9 = @(42400, 40900)
    10 = @(49100, 44400)
    12 = @(7400, 7700)
    13 = @(0, 0)
    14 = @(2700, 4400)
    15 = @(5200, 4100)
    17 = @(87500, 85000)
    18 = @(4000, 300)
    20 = @(100, 200)
    21 = @(10800, 6000)
    22 = @(1300, 900)
    23 = @(2700, -500)
    24 = @(600, 0)
    25 = @(0, 0)
    26 = @(2100, -500)
    27 = @(2100, -500)
    28 = @(0, 0)
    29 = @("NA", "NA")
    30 = @(0, 0)
    31 = @(0, 0)
    32 = @(-100, -200)
    33 = @(2100, -500)
    34 = @(0, 0)
    35 = @(2100, -500)
    38 = @(43434, 43343)
    41 = @(42800, 24800)
    42 = @(1400, 1400)
    43 = @(43400, 40200)
    44 = @(50600, 49700)
    45 = @(4800, 4900)
    46 = @(143000, 120900)
    47 = @(0, 0)
    48 = @(41900, 42200)
    49 = @(595600, 527200)
    50 = @(0, 0)
    51 = @(0, 0)
    52 = @(3500, 3600)
    53 = @(0, 0)
    54 = @(784000, 693900)
    57 = @(19400, 16900)
    58 = @(5000, 5000)
    59 = @(25300, 23900)
    60 = @(49700, 45800)
    61 = @(139300, 85400)
    62 = @(45400, 17700)
    63 = @(0, 0)
    64 = @(0, 0)
    65 = @(0, 0)
    66 = @(234400, 149000)
    68 = @(0, 0)
    69 = @(0, 0)
    70 = @(0, 0)
    71 = @(0, 0)
    72 = @(6800, 4700)
    73 = @(0, 0)
    74 = @(0, 0)
    75 = @(0, 0)
    76 = @(549600, 544900)
    77 = @(0, 0)
    80 = @(43434, 43343)
    81 = @(2100, -500)
    83 = @(6700, 5600)
    84 = @(0, 0)
    85 = @(0, 0)
    86 = @(0, 0)
    87 = @(0, 0)
    88 = @(0, 0)
    89 = @(13000, -8900)
    91 = @(-700, -700)
    92 = @(0, 0)
    93 = @(0, 0)
    94 = @(-48700, -37700)
    96 = @(0, 0)
    97 = @(0, 0)
    98 = @(0, 0)
    99 = @(0, 0)
    100 = @(53900, -2600)
    101 = @(-100, -100)
    102 = @(18100, -49300)
}
foreach ($r in $newData.Keys) {
    $vals = $newData[$r]
    $ws.Cells.Item([int]$r, 4).Value2 = $vals[0]
    $ws.Cells.Item([int]$r, 5).Value2 = $vals[1]
}
